$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Update revised weekly figures (column G, Waargenomen) for weeks already
#    present. The formula in column I ("Oversterfte" = G - H) recalculates
#    automatically.
$ws.Range("G22").Value = 2673
$ws.Range("G23").Value = 2668
$ws.Range("G30").Value = 2719
$ws.Range("G32").Value = 2997
$ws.Range("G34").Value = 3216
$ws.Range("G35").Value = 3445
$ws.Range("G36").Value = 3675
$ws.Range("G38").Value = 3560
$ws.Range("G39").Value = 3317
$ws.Range("G40").Value = 3388
$ws.Range("G41").Value = 3494

# 2) Make room for 3 new weeks (50, 51, 52) by inserting rows just above the
#    totals row (currently row 43), shifting the totals row down to row 46.
$ws.Range("A43:I45").Insert()

# 3) Fill in week 50 (full data row, continuing the "Oversterfte" formula).
$ws.Range("F42").Value = 50
$ws.Range("G42").Value = 3571
$ws.Range("H42").Value = 3100
$ws.Range("I42").Formula = "=G42-H42"

# 4) Weeks 51 and 52 only have the week number filled in so far.
$ws.Range("F43").Value = 51
$ws.Range("F44").Value = 52

# 5) Restore the view state recorded after the edit.
$excel.ActiveWindow.ScrollRow = 22
$ws.Range("I43").Select() | Out-Null
